# Excel template bug fix
# - Clear the "Obs_relatorio" (column H) validation messages on sheet "Bico"
#   for rows 2-7 (the template no longer prints the success message there).
# - On sheet "Tanque", fix the "Fechamento" value for row 2 and recompute the
#   validation messages (column F) to report the real divergences between the
#   SPED total and the report values for rows 2-4.

$wb = $excel.ActiveWorkbook

# --- Sheet "Bico" ---------------------------------------------------------
$bico = $wb.Worksheets.Item("Bico")
foreach ($r in 2..7) {
    $bico.Cells.Item($r, 8).Value = ""
}

# --- Sheet "Tanque" --------------------------------------------------------
$tanque = $wb.Worksheets.Item("Tanque")

# Fechamento value bugfix for row 2
$tanque.Cells.Item(2, 4).Value = 200

# Updated divergence messages (column F) for rows 2-4
$tanque.Cells.Item(2, 6).Value = "Divergência entre o SPED(4056,00) e o relatório(200,00)!"
$tanque.Cells.Item(3, 6).Value = "Divergência entre o SPED(4056,00) e o relatório(4157,00)!"
$tanque.Cells.Item(4, 6).Value = "Divergência entre o SPED(4056,00) e o relatório(6333,00)!"
